# Update "Name of Algo" - change values in column B for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 6.439
    13 = 5.302999999999999
    16 = 5.294
    18 = 5.225000000000001
    20 = 6.24
    26 = 6.516
    27 = 6.161
    29 = 5.607000000000001
    35 = 8.317
    36 = 7.972
    45 = 5.627
    55 = 4.711
    57 = 5.147
    69 = 5.013
    76 = 6.157000000000001
    78 = 7.811
    82 = 5.286
    83 = 5.667
    93 = 5.377
    97 = 5.632
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
